$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-12-03 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-04 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("98-50=48", $true, $false, $false, $false, $false, $true, 1, $false, "51-44=7", 2) | Out-Null
$d.Content.Find.Execute("58+24=82", $true, $false, $false, $false, $false, $true, 1, $false, "81-76=5", 2) | Out-Null
$d.Content.Find.Execute("7+30=37", $true, $false, $false, $false, $false, $true, 1, $false, "53-5=48", 2) | Out-Null
$d.Content.Find.Execute("18+29=47", $true, $false, $false, $false, $false, $true, 1, $false, "87-85=2", 2) | Out-Null
$d.Content.Find.Execute("67-24=43", $true, $false, $false, $false, $false, $true, 1, $false, "85-61=24", 2) | Out-Null
$d.Content.Find.Execute("42+8=50", $true, $false, $false, $false, $false, $true, 1, $false, "71+24=95", 2) | Out-Null
$d.Content.Find.Execute("23+34=57", $true, $false, $false, $false, $false, $true, 1, $false, "42-29=13", 2) | Out-Null
$d.Content.Find.Execute("17+41=58", $true, $false, $false, $false, $false, $true, 1, $false, "72+6=78", 2) | Out-Null
$d.Content.Find.Execute("92-9=83", $true, $false, $false, $false, $false, $true, 1, $false, "83+12=95", 2) | Out-Null
$d.Content.Find.Execute("30+63=93", $true, $false, $false, $false, $false, $true, 1, $false, "93-31=62", 2) | Out-Null
$d.Content.Find.Execute("91-64=27", $true, $false, $false, $false, $false, $true, 1, $false, "36+59=95", 2) | Out-Null
$d.Content.Find.Execute("0+97=97", $true, $false, $false, $false, $false, $true, 1, $false, "33-8=25", 2) | Out-Null
$d.Content.Find.Execute("25-6=19", $true, $false, $false, $false, $false, $true, 1, $false, "92-26=66", 2) | Out-Null
$d.Content.Find.Execute("77-52=25", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=71", 2) | Out-Null
$d.Content.Find.Execute("53-26=27", $true, $false, $false, $false, $false, $true, 1, $false, "3+95=98", 2) | Out-Null
$d.Content.Find.Execute("4+77=81", $true, $false, $false, $false, $false, $true, 1, $false, "79-32=47", 2) | Out-Null
$d.Content.Find.Execute("91-79=12", $true, $false, $false, $false, $false, $true, 1, $false, "99-36=63", 2) | Out-Null
$d.Content.Find.Execute("11+23=34", $true, $false, $false, $false, $false, $true, 1, $false, "15-8=7", 2) | Out-Null
$d.Content.Find.Execute("53-40=13", $true, $false, $false, $false, $false, $true, 1, $false, "64+23=87", 2) | Out-Null
$d.Content.Find.Execute("44+39=83", $true, $false, $false, $false, $false, $true, 1, $false, "58-37=21", 2) | Out-Null
$d.Content.Find.Execute("19-14=5", $true, $false, $false, $false, $false, $true, 1, $false, "69+29=98", 2) | Out-Null
$d.Content.Find.Execute("31+24=55", $true, $false, $false, $false, $false, $true, 1, $false, "48+24=72", 2) | Out-Null
$d.Content.Find.Execute("66-20=46", $true, $false, $false, $false, $false, $true, 1, $false, "7+19=26", 2) | Out-Null
$d.Content.Find.Execute("45+6=51", $true, $false, $false, $false, $false, $true, 1, $false, "77-74=3", 2) | Out-Null
$d.Content.Find.Execute("66+32=98", $true, $false, $false, $false, $false, $true, 1, $false, "51+15=66", 2) | Out-Null
$d.Content.Find.Execute("9+81=90", $true, $false, $false, $false, $false, $true, 1, $false, "69-54=15", 2) | Out-Null
$d.Content.Find.Execute("4+48=52", $true, $false, $false, $false, $false, $true, 1, $false, "50+29=79", 2) | Out-Null
$d.Content.Find.Execute("13+60=73", $true, $false, $false, $false, $false, $true, 1, $false, "57+9=66", 2) | Out-Null
$d.Content.Find.Execute("28+47=75", $true, $false, $false, $false, $false, $true, 1, $false, "91-89=2", 2) | Out-Null
$d.Content.Find.Execute("17+8=25", $true, $false, $false, $false, $false, $true, 1, $false, "42+29=71", 2) | Out-Null
$d.Content.Find.Execute("14-11=3", $true, $false, $false, $false, $false, $true, 1, $false, "62-5=57", 2) | Out-Null
$d.Content.Find.Execute("51-6=45", $true, $false, $false, $false, $false, $true, 1, $false, "5+43=48", 2) | Out-Null
$d.Content.Find.Execute("11+2=13", $true, $false, $false, $false, $false, $true, 1, $false, "0+18=18", 2) | Out-Null
$d.Content.Find.Execute("8+89=97", $true, $false, $false, $false, $false, $true, 1, $false, "86-71=15", 2) | Out-Null
$d.Content.Find.Execute("14+56=70", $true, $false, $false, $false, $false, $true, 1, $false, "69-18=51", 2) | Out-Null
$d.Content.Find.Execute("91+3=94", $true, $false, $false, $false, $false, $true, 1, $false, "43-25=18", 2) | Out-Null
$d.Content.Find.Execute("58-18=40", $true, $false, $false, $false, $false, $true, 1, $false, "48+25=73", 2) | Out-Null
$d.Content.Find.Execute("66-36=30", $true, $false, $false, $false, $false, $true, 1, $false, "12-4=8", 2) | Out-Null
$d.Content.Find.Execute("8+87=95", $true, $false, $false, $false, $false, $true, 1, $false, "32-6=26", 2) | Out-Null
$d.Content.Find.Execute("88-45=43", $true, $false, $false, $false, $false, $true, 1, $false, "94-49=45", 2) | Out-Null
$d.Content.Find.Execute("86-60=26", $true, $false, $false, $false, $false, $true, 1, $false, "39+26=65", 2) | Out-Null
$d.Content.Find.Execute("20+50=70", $true, $false, $false, $false, $false, $true, 1, $false, "28+31=59", 2) | Out-Null
$d.Content.Find.Execute("32+6=38", $true, $false, $false, $false, $false, $true, 1, $false, "0+5=5", 2) | Out-Null
$d.Content.Find.Execute("73+11=84", $true, $false, $false, $false, $false, $true, 1, $false, "29+25=54", 2) | Out-Null
$d.Content.Find.Execute("92-22=70", $true, $false, $false, $false, $false, $true, 1, $false, "90-8=82", 2) | Out-Null
$d.Content.Find.Execute("59-57=2", $true, $false, $false, $false, $false, $true, 1, $false, "59+22=81", 2) | Out-Null
$d.Content.Find.Execute("53-42=11", $true, $false, $false, $false, $false, $true, 1, $false, "44+28=72", 2) | Out-Null
$d.Content.Find.Execute("16+15=31", $true, $false, $false, $false, $false, $true, 1, $false, "29+16=45", 2) | Out-Null
$d.Content.Find.Execute("2+69=71", $true, $false, $false, $false, $false, $true, 1, $false, "79-75=4", 2) | Out-Null
$d.Content.Find.Execute("18+42=60", $true, $false, $false, $false, $false, $true, 1, $false, "6+80=86", 2) | Out-Null
$d.Content.Find.Execute("88-30=58", $true, $false, $false, $false, $false, $true, 1, $false, "47+32=79", 2) | Out-Null
$d.Content.Find.Execute("50-13=37", $true, $false, $false, $false, $false, $true, 1, $false, "95-81=14", 2) | Out-Null
$d.Content.Find.Execute("85-83=2", $true, $false, $false, $false, $false, $true, 1, $false, "70-29=41", 2) | Out-Null
$d.Content.Find.Execute("36+11=47", $true, $false, $false, $false, $false, $true, 1, $false, "92-19=73", 2) | Out-Null
$d.Content.Find.Execute("38+35=73", $true, $false, $false, $false, $false, $true, 1, $false, "6+64=70", 2) | Out-Null
$d.Content.Find.Execute("53+18=71", $true, $false, $false, $false, $false, $true, 1, $false, "89+1=90", 2) | Out-Null
$d.Content.Find.Execute("23+8=31", $true, $false, $false, $false, $false, $true, 1, $false, "56-55=1", 2) | Out-Null
$d.Content.Find.Execute("61+38=99", $true, $false, $false, $false, $false, $true, 1, $false, "10+49=59", 2) | Out-Null
$d.Content.Find.Execute("47+5=52", $true, $false, $false, $false, $false, $true, 1, $false, "59-1=58", 2) | Out-Null
$d.Content.Find.Execute("30-26=4", $true, $false, $false, $false, $false, $true, 1, $false, "41+55=96", 2) | Out-Null
$d.Content.Find.Execute("46-37=9", $true, $false, $false, $false, $false, $true, 1, $false, "30-7=23", 2) | Out-Null
$d.Content.Find.Execute("87-35=52", $true, $false, $false, $false, $false, $true, 1, $false, "89-13=76", 2) | Out-Null
$d.Content.Find.Execute("92-77=15", $true, $false, $false, $false, $false, $true, 1, $false, "14+35=49", 2) | Out-Null
$d.Content.Find.Execute("29+64=93", $true, $false, $false, $false, $false, $true, 1, $false, "93-6=87", 2) | Out-Null
$d.Content.Find.Execute("17-11=6", $true, $false, $false, $false, $false, $true, 1, $false, "40+49=89", 2) | Out-Null
$d.Content.Find.Execute("10+21=31", $true, $false, $false, $false, $false, $true, 1, $false, "85-77=8", 2) | Out-Null
$d.Content.Find.Execute("15+80=95", $true, $false, $false, $false, $false, $true, 1, $false, "70-64=6", 2) | Out-Null
$d.Content.Find.Execute("76+0=76", $true, $false, $false, $false, $false, $true, 1, $false, "64+19=83", 2) | Out-Null
$d.Content.Find.Execute("54+10=64", $true, $false, $false, $false, $false, $true, 1, $false, "86-4=82", 2) | Out-Null
$d.Content.Find.Execute("33+13=46", $true, $false, $false, $false, $false, $true, 1, $false, "26+49=75", 2) | Out-Null
$d.Content.Find.Execute("86-55=31", $true, $false, $false, $false, $false, $true, 1, $false, "28+56=84", 2) | Out-Null
$d.Content.Find.Execute("12+86=98", $true, $false, $false, $false, $false, $true, 1, $false, "99-55=44", 2) | Out-Null
$d.Content.Find.Execute("87-87=0", $true, $false, $false, $false, $false, $true, 1, $false, "92-49=43", 2) | Out-Null
$d.Content.Find.Execute("39+50=89", $true, $false, $false, $false, $false, $true, 1, $false, "63-55=8", 2) | Out-Null
$d.Content.Find.Execute("4-4=0", $true, $false, $false, $false, $false, $true, 1, $false, "81-63=18", 2) | Out-Null
$d.Content.Find.Execute("54-33=21", $true, $false, $false, $false, $false, $true, 1, $false, "80-19=61", 2) | Out-Null
$d.Content.Find.Execute("89-74=15", $true, $false, $false, $false, $false, $true, 1, $false, "1+25=26", 2) | Out-Null
$d.Content.Find.Execute("98-98=0", $true, $false, $false, $false, $false, $true, 1, $false, "90-81=9", 2) | Out-Null
$d.Content.Find.Execute("71+13=84", $true, $false, $false, $false, $false, $true, 1, $false, "7+83=90", 2) | Out-Null
$d.Content.Find.Execute("77+7=84", $true, $false, $false, $false, $false, $true, 1, $false, "27+45=72", 2) | Out-Null
$d.Content.Find.Execute("25-17=8", $true, $false, $false, $false, $false, $true, 1, $false, "83-79=4", 2) | Out-Null
$d.Content.Find.Execute("71-21=50", $true, $false, $false, $false, $false, $true, 1, $false, "15+48=63", 2) | Out-Null
$d.Content.Find.Execute("25+3=28", $true, $false, $false, $false, $false, $true, 1, $false, "77-56=21", 2) | Out-Null
$d.Content.Find.Execute("30-13=17", $true, $false, $false, $false, $false, $true, 1, $false, "73-54=19", 2) | Out-Null
$d.Content.Find.Execute("2+82=84", $true, $false, $false, $false, $false, $true, 1, $false, "24-13=11", 2) | Out-Null
$d.Content.Find.Execute("17+60=77", $true, $false, $false, $false, $false, $true, 1, $false, "31+37=68", 2) | Out-Null
$d.Content.Find.Execute("69-49=20", $true, $false, $false, $false, $false, $true, 1, $false, "27-0=27", 2) | Out-Null
$d.Content.Find.Execute("41-28=13", $true, $false, $false, $false, $false, $true, 1, $false, "48-34=14", 2) | Out-Null
$d.Content.Find.Execute("29-10=19", $true, $false, $false, $false, $false, $true, 1, $false, "59-47=12", 2) | Out-Null
$d.Content.Find.Execute("75-16=59", $true, $false, $false, $false, $false, $true, 1, $false, "80-49=31", 2) | Out-Null
$d.Content.Find.Execute("97-0=97", $true, $false, $false, $false, $false, $true, 1, $false, "45-43=2", 2) | Out-Null
$d.Content.Find.Execute("35+55=90", $true, $false, $false, $false, $false, $true, 1, $false, "29-7=22", 2) | Out-Null
$d.Content.Find.Execute("66+18=84", $true, $false, $false, $false, $false, $true, 1, $false, "88+9=97", 2) | Out-Null
$d.Content.Find.Execute("86-24=62", $true, $false, $false, $false, $false, $true, 1, $false, "87-12=75", 2) | Out-Null
$d.Content.Find.Execute("30+52=82", $true, $false, $false, $false, $false, $true, 1, $false, "0+73=73", 2) | Out-Null
$d.Content.Find.Execute("86-17=69", $true, $false, $false, $false, $false, $true, 1, $false, "42+54=96", 2) | Out-Null
$d.Content.Find.Execute("36-6=30", $true, $false, $false, $false, $false, $true, 1, $false, "67+4=71", 2) | Out-Null
$d.Content.Find.Execute("93-89=4", $true, $false, $false, $false, $false, $true, 1, $false, "81-35=46", 2) | Out-Null
$d.Content.Find.Execute("42+26=68", $true, $false, $false, $false, $false, $true, 1, $false, "98-32=66", 2) | Out-Null
$d.Content.Find.Execute("13+7=20", $true, $false, $false, $false, $false, $true, 1, $false, "57+29=86", 2) | Out-Null
